# Fix the misspelled Portuguese word "urgencia" -> "urgência" and drop the
# spell-check flags (w:proofErr spellStart/spellEnd) that wrapped it, since
# the word is now spelled correctly.
#
# A plain Find/Replace limited to the misspelled run's own text leaves the
# surrounding proofErr markers untouched. So we first do a replace whose
# matched range spans across both run boundaries (absorbing the proofErr
# markers along with it, which merges the three runs into one) and then
# "touch" the formatting of just the corrected word so Word re-splits that
# merged run back into three separate runs - matching the original run
# layout, but now proofErr-free.

$d = $word.ActiveDocument

# Step 1: replace across the run boundaries on both sides of "urgencia" so
# the proofErr spellStart/spellEnd pair bracketing it gets swallowed too.
$fix = $d.Content
$fix.Find.Execute("serviço de urgencia  com", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "serviço de urgência  com", 2)

# Step 2: re-split the now-merged run so "urgência" is its own run again
# (matching the original document's run structure) by toggling a character
# format on just that word.
$rng = $d.Content
$rng.Find.Execute("urgência", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$rng.Font.Bold = $true
$rng.Font.Bold = $false
